$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, border, centered) from an existing header cell (AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the team record values for every data row (2-43)
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 88   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 74   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 1    # AF - Ties
}
